$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.04
$wsSummary.Range("B4").Value = 0.03
$wsSummary.Range("B5").Value = 0.05
$wsSummary.Range("B6").Value = 11
$wsSummary.Range("B7").Value = 5
$wsSummary.Range("B9").Value = 45.45

# --- Strategy Status sheet ---
$wsStrategy = $wb.Worksheets.Item("Strategy Status")
$wsStrategy.Range("C4").Value = 100.04
$wsStrategy.Range("D4").Value = 11
$wsStrategy.Range("E4").Value = 0.03
$wsStrategy.Range("F4").Value = 0.04
$wsStrategy.Range("G4").Value = 45.45

# --- Helper to append the new trade row #11 (worksheet row 12) ---
function Add-TradeRow($ws) {
    $ws.Cells.Item(12, 1).Value = 11
    $ws.Cells.Item(12, 2).NumberFormat = "@"
    $ws.Cells.Item(12, 2).Value = "2026-02-17"
    $ws.Cells.Item(12, 3).NumberFormat = "@"
    $ws.Cells.Item(12, 3).Value = "12:27:37"
    $ws.Cells.Item(12, 4).Value = "MarketMaking"
    $ws.Cells.Item(12, 5).Value = "UP"
    $ws.Cells.Item(12, 6).Value = 0.08
    $ws.Cells.Item(12, 7).Value = 0.093468
    $ws.Cells.Item(12, 8).Value = "CLOSED"
    $ws.Cells.Item(12, 9).Value = 16.8355
    $ws.Cells.Item(12, 10).Value = 0.01
    $ws.Cells.Item(12, 11).Value = 100.04
    $ws.Cells.Item(12, 12).Value = 0
    $ws.Cells.Item(12, 13).Value = 0
    $ws.Cells.Item(12, 14).Value = 0.6
    $ws.Cells.Item(12, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(12, 16).Value = "early_exit"
    $ws.Cells.Item(12, 17).Value = 0.13
}

# --- All Trades sheet ---
$wsAllTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAllTrades

# --- MarketMaking sheet ---
$wsMarketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow $wsMarketMaking
